$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47 currently holds "Lugo" data, row 48 currently holds "Almeria" data.
# The update swaps the two whole rows (province name + all four numeric
# columns) so that "Almeria" ends up on row 47 and "Lugo" on row 48.
$lugoName = $ws.Range("A47").Text
$lugoB    = $ws.Range("B47").Value2
$lugoC    = $ws.Range("C47").Value2
$lugoD    = $ws.Range("D47").Value2
$lugoE    = $ws.Range("E47").Value2

$almeriaName = $ws.Range("A48").Text
$almeriaB    = $ws.Range("B48").Value2
$almeriaC    = $ws.Range("C48").Value2
$almeriaD    = $ws.Range("D48").Value2
$almeriaE    = $ws.Range("E48").Value2

$ws.Range("A47").Value = $almeriaName
$ws.Range("B47").Value = $almeriaB
$ws.Range("C47").Value = $almeriaC
$ws.Range("D47").Value = $almeriaD
$ws.Range("E47").Value = $almeriaE

$ws.Range("A48").Value = $lugoName
$ws.Range("B48").Value = $lugoB
$ws.Range("C48").Value = $lugoC
$ws.Range("D48").Value = $lugoD
$ws.Range("E48").Value = $lugoE

# Update the "last updated" timestamp string in cell A1.
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 23:16"
